# Fix the tillage/application-date placeholders on the "scenario" sheet.
# The template had every year row pointing at the same "01/01/YYYY" shared
# string for tillage date (E), N-application date (J) and a third date (L).
# Correct values per row:
#   E{row} = "10/01/<year>"   (tillage_date)
#   J{row} = "05/01/<year>"   (n_application_date)
#   L{row} = "05/09/<year>"   (the secondary n_application date column)
# Also B7 ("crop_scenario_name"-ish lookup cell) was left at "Current" and
# should be "no till".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario")

$dataRows = @(16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,37,38,39,40,41,42,43,44,45,46,48,49,50,51,52,53,54,55,56,57)

foreach ($r in $dataRows) {
    $year = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 5).Value = '"10/01/' + $year + '"'
    $ws.Cells.Item($r, 10).Value = '"05/01/' + $year + '"'
    $ws.Cells.Item($r, 12).Value = '"05/09/' + $year + '"'
}

$ws.Range("B7").Value = "no till"

# Restore the view: scroll back to the top and leave a single-cell selection
# on M50 instead of the old D48:D57 block.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("M50").Select() | Out-Null
